$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New note cells alongside rows 2 and 3 (column Q)
$ws.Range("Q2").Value = "to 1/4/2016"
$ws.Range("Q3").Value = "!test on open price similar to all ticks"

# New assessment row (row 5) - Hc1/ac1 test, not beneficial
$ws.Range("A5").Value = "test Hc1ac1"
$ws.Range("F5").Value = 30
$ws.Range("G5").Value = 30
$ws.Range("N5").Value = 90
$ws.Range("O5").Value = 750
$ws.Range("P5").Value = 1.3

# Reset view: scroll back to top-left and move selection to O7
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("O7").Select()
